$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CMF")

# Rename the "Index" column header to "i" (also renames the "testdata" table's
# first column, since the table header row is backed by this same cell)
$ws.Range("A1").Value = "i"

# The Index column used to be 1-based (1, 2, 3, ...); re-base it to 0-based
# (0, 1, 2, ...) by decrementing every data row's value by 1.
$lastRow = $ws.UsedRange.Rows.Count()
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    $cell.Value = $current - 1
}

# Narrow column A now that the header text is shorter ("i" vs "Index");
# ColumnWidth=3.1666666667 (character units) serializes to a stored width of 4.
$ws.Columns.Item(1).ColumnWidth = 3.1666666667
